$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each of these rows corresponds to a newly-included 2020-08-24 data point
# for the Fonds de solidarite volet 2 (regional / categorie juridique) dataset.
# "nombre_aides" (column C) and "montant_total" (column D) are updated to the
# latest published counts/amounts. Values are written as text (matching the
# inlineStr / shared-string format already used throughout the sheet), and the
# temporary text number-format is reverted afterwards so no cell style changes
# are left behind.

$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "695"
$ws.Range("C2").Style = "Normal"

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "1597687.79"
$ws.Range("D2").Style = "Normal"

$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "1018"
$ws.Range("C4").Style = "Normal"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "3609453.47"
$ws.Range("D4").Style = "Normal"

$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "655"
$ws.Range("C6").Style = "Normal"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2109277.78"
$ws.Range("D6").Style = "Normal"

$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "221"
$ws.Range("C14").Style = "Normal"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "598362.00"
$ws.Range("D14").Style = "Normal"

$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "497"
$ws.Range("C16").Style = "Normal"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1836574.13"
$ws.Range("D16").Style = "Normal"

$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = "336"
$ws.Range("C21").Style = "Normal"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1185741.00"
$ws.Range("D21").Style = "Normal"

$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "161"
$ws.Range("C22").Style = "Normal"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "468812.39"
$ws.Range("D22").Style = "Normal"

$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "282"
$ws.Range("C28").Style = "Normal"

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "728137.45"
$ws.Range("D28").Style = "Normal"

$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "159"
$ws.Range("C40").Style = "Normal"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "431017.22"
$ws.Range("D40").Style = "Normal"

$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "86"
$ws.Range("C41").Style = "Normal"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "415409.98"
$ws.Range("D41").Style = "Normal"

$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "127"
$ws.Range("C42").Style = "Normal"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "520972.99"
$ws.Range("D42").Style = "Normal"

$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "420"
$ws.Range("C48").Style = "Normal"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1451998.40"
$ws.Range("D48").Style = "Normal"

$ws.Range("C76").NumberFormat = "@"
$ws.Range("C76").Value = "928"
$ws.Range("C76").Style = "Normal"

$ws.Range("D76").NumberFormat = "@"
$ws.Range("D76").Value = "3245240.26"
$ws.Range("D76").Style = "Normal"

$ws.Range("C77").NumberFormat = "@"
$ws.Range("C77").Value = "523"
$ws.Range("C77").Style = "Normal"

$ws.Range("D77").NumberFormat = "@"
$ws.Range("D77").Value = "1724125.47"
$ws.Range("D77").Style = "Normal"

$ws.Range("C96").NumberFormat = "@"
$ws.Range("C96").Value = "1017"
$ws.Range("C96").Style = "Normal"

$ws.Range("D96").NumberFormat = "@"
$ws.Range("D96").Value = "3122636.31"
$ws.Range("D96").Style = "Normal"

